$wb = $excel.ActiveWorkbook

# Map of F-column (想去人数) updates shared by both "展览" and "全部类型" sheets.
$fUpdates = @{
    2  = 1110
    3  = 427
    4  = 1520
    5  = 8801
    7  = 496
    9  = 298
    11 = 23
    12 = 18
    13 = 3662
    15 = 370
    16 = 85
    17 = 2701
    20 = 316
    21 = 215
    22 = 2446
}

# Sheet "展览" - row 23 holds the last record (71 -> 77 / 88 -> 128)
$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $fUpdates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $fUpdates[$row]
}
$ws1.Cells.Item(23, 6).Value = 77
$ws1.Cells.Item(23, 7).Value = 128

# Sheet "全部类型" - same F-column updates, but the last record is row 24 (71 -> 77 / 88 -> 128)
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $fUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $fUpdates[$row]
}
$ws4.Cells.Item(24, 6).Value = 77
$ws4.Cells.Item(24, 7).Value = 128
